# Update last_edited_time (column D) for rows 2-22 from 15:02 to 16:53 UTC
# as part of the "code for KPI update" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2024-07-25T16:53:00.000Z"

for ($row = 2; $row -le 22; $row++) {
    $ws.Cells.Item($row, 4).Value = $newTimestamp
}
